$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit message: "add the NA's under duplicate_image_filename"
# Column E is "duplicate_image_filename"; rows 2-21 hold the practice (p1-p4)
# and generic/unique_video/unique_audio stimulus rows that previously had
# no value in that column. Fill them with "NA".
$ws.Range("E2:E21").Value = "NA"
